$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value2 = '67.848.01'
$cell.Style = "Normal"
$ws.Range("E2").Value2 = '  -0.95%  '
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value2 = '3.849.16'
$cell.Style = "Normal"
$ws.Range("E3").Value2 = '  -1.47%  '
$ws.Range("E4").Value2 = '  +0.04%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value2 = '598.23'
$cell.Style = "Normal"
$ws.Range("E5").Value2 = '  -0.68%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value2 = '166.06'
$cell.Style = "Normal"
$ws.Range("E6").Value2 = '  +1.13%  '
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value2 = '3.848.84'
$cell.Style = "Normal"
$ws.Range("E7").Value2 = '  -1.42%  '
$ws.Range("E8").Value2 = '  -0.07%  '
$ws.Range("E9").Value2 = '  +0.01%  '
$ws.Range("E10").Value2 = '  -0.43%  '
$ws.Range("E11").Value2 = '  -0.81%  '
$ws.Range("E12").Value2 = '  -0.45%  '
$ws.Range("E13").Value2 = '  +0.65%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value2 = '36.74'
$cell.Style = "Normal"
$ws.Range("E14").Value2 = '  +0.12%  '
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value2 = '4.496.04'
$cell.Style = "Normal"
$ws.Range("E15").Value2 = '  -1.38%  '
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value2 = '3.864.02'
$cell.Style = "Normal"
$ws.Range("E16").Value2 = '  -1.72%  '
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value2 = '67.902.85'
$cell.Style = "Normal"
$ws.Range("E17").Value2 = '  -1.11%  '
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value2 = '18.07'
$cell.Style = "Normal"
$ws.Range("E18").Value2 = '  +6.48%  '
$ws.Range("E19").Value2 = '  -0.61%  '
$ws.Range("E20").Value2 = '  -1.42%  '
$ws.Range("E21").Value2 = '  -2.83%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value2 = '462.88'
$cell.Style = "Normal"
$ws.Range("E22").Value2 = '  -4.16%  '
$ws.Range("E23").Value2 = '  +1.69%  '
$ws.Range("E24").Value2 = '  -4.00%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value2 = '83.14'
$cell.Style = "Normal"
$ws.Range("E25").Value2 = '  -1.35%  '
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value2 = '2.23'
$cell.Style = "Normal"
$ws.Range("E26").Value2 = '  +0.39%  '
$ws.Range("E27").Value2 = '  +1.03%  '
$ws.Range("E28").Value2 = '  -0.01%  '
$ws.Range("E29").Value2 = '  -0.87%  '
$ws.Range("E30").Value2 = '  +0.47%  '
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value2 = '4.002.33'
$cell.Style = "Normal"
$ws.Range("E31").Value2 = '  -1.32%  '
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value2 = '7.72'
$cell.Style = "Normal"
$ws.Range("E32").Value2 = '  -1.14%  '
$ws.Range("E33").Value2 = '  -2.00%  '
$ws.Range("E34").Value2 = '  -2.77%  '
$ws.Range("B35").Value2 = 'Aptos'
$ws.Range("C35").Value2 = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value2 = '9.30'
$cell.Style = "Normal"
$ws.Range("E35").Value2 = '  +2.34%  '
$ws.Range("B36").Value2 = 'RenzoRestakedETH'
$ws.Range("C36").Value2 = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value2 = '3.828.43'
$cell.Style = "Normal"
$ws.Range("E36").Value2 = '  -0.54%  '
$ws.Range("B37").Value2 = 'Hedera'
$ws.Range("C37").Value2 = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value2 = '0.104'
$cell.Style = "Normal"
$ws.Range("E37").Value2 = '  -2.10%  '
$ws.Range("B38").Value2 = 'Mantle'
$ws.Range("C38").Value2 = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value2 = '1.02'
$cell.Style = "Normal"
$ws.Range("E38").Value2 = '  -1.68%  '
$ws.Range("B39").Value2 = 'Kaspa'
$ws.Range("C39").Value2 = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value2 = '0.139'
$cell.Style = "Normal"
$ws.Range("E39").Value2 = '  -0.02%  '
$ws.Range("B40").Value2 = 'Filecoin'
$ws.Range("C40").Value2 = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value2 = '5.88'
$cell.Style = "Normal"
$ws.Range("E40").Value2 = '  +0.49%  '
$ws.Range("B41").Value2 = 'dogwifhat'
$ws.Range("C41").Value2 = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value2 = '3.27'
$cell.Style = "Normal"
$ws.Range("E41").Value2 = '  +7.17%  '
$ws.Range("B42").Value2 = 'FirstDigitalUSD'
$ws.Range("C42").Value2 = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value2 = '1.00'
$cell.Style = "Normal"
$ws.Range("E42").Value2 = '  +0.10%  '
$ws.Range("B43").Value2 = 'TheGraph'
$ws.Range("C43").Value2 = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value2 = '0.312'
$cell.Style = "Normal"
$ws.Range("E43").Value2 = '  -1.29%  '
$ws.Range("B44").Value2 = 'Bittensor'
$ws.Range("C44").Value2 = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value2 = '426.81'
$cell.Style = "Normal"
$ws.Range("E44").Value2 = '  -0.99%  '
$ws.Range("B45").Value2 = 'Stacks'
$ws.Range("C45").Value2 = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value2 = '1.97'
$cell.Style = "Normal"
$ws.Range("E45").Value2 = '  +0.09%  '
$ws.Range("B46").Value2 = 'USDe'
$ws.Range("C46").Value2 = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value2 = '1.00'
$cell.Style = "Normal"
$ws.Range("E46").Value2 = '  -0.02%  '
$ws.Range("B47").Value2 = 'OKB'
$ws.Range("C47").Value2 = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value2 = '47.13'
$cell.Style = "Normal"
$ws.Range("E47").Value2 = '  -2.78%  '
$ws.Range("B48").Value2 = 'Cosmos'
$ws.Range("C48").Value2 = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value2 = '8.48'
$cell.Style = "Normal"
$ws.Range("E48").Value2 = '  +0.97%  '
$ws.Range("B49").Value2 = 'EnergySwap'
$ws.Range("C49").Value2 = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value2 = '26.67'
$cell.Style = "Normal"
$ws.Range("E49").Value2 = '  +0.96%  '
$ws.Range("B50").Value2 = 'FLOKI'
$ws.Range("C50").Value2 = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value2 = '0.000274'
$cell.Style = "Normal"
$ws.Range("E50").Value2 = '  +4.04%  '
$ws.Range("B51").Value2 = 'Monero'
$ws.Range("C51").Value2 = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value2 = '143.54'
$cell.Style = "Normal"
$ws.Range("E51").Value2 = '  +1.20%  '
